$wb = $excel.ActiveWorkbook

# Helper: set a cell so it is stored as TEXT (shared string), matching the
# original workbook's convention of keeping numeric-looking values as text,
# while avoiding leaving a lingering cell-level style reference.
function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# ---------------------------------------------------------------
# Sheet "Restricciones_del_follower" (3rd sheet) - J_0_L0_v / J_0_LP_v / J_Ne_L0_v rows
# ---------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item(3)

Set-TextCell $wsFollower "A2" "-12.772972972972951 - 2x_1 + 1.135135135135135y_1 - 0.810810810810811y_2"
Set-TextCell $wsFollower "B2" "15.272972972972951"
Set-TextCell $wsFollower "D2" "0.77"
Set-TextCell $wsFollower "E2" "3.3000000000000003"
Set-TextCell $wsFollower "F2" "0"

Set-TextCell $wsFollower "A3" "173.5641891891892 + x_1 - 3x_2 - 1.5540540540540542y_1 - 1.1756756756756754y_2"
Set-TextCell $wsFollower "B3" "-175.5641891891892"
Set-TextCell $wsFollower "D3" "0.46"
Set-TextCell $wsFollower "E3" "-1.6"
Set-TextCell $wsFollower "F3" "-2.4"

Set-TextCell $wsFollower "A4" "57.75074849075653 - 0.5640713547106777y_1 + 0.6103001034050511y_2"
Set-TextCell $wsFollower "B4" "-57.75074849075653"
Set-TextCell $wsFollower "D4" "0.41"
Set-TextCell $wsFollower "E4" "4.4"
Set-TextCell $wsFollower "F4" "0"

Set-TextCell $wsFollower "A5" "-368.83581081081087 + 3.4459459459459465y_1 + 3.8243243243243246y_2"
Set-TextCell $wsFollower "B5" "368.28581081081086"
Set-TextCell $wsFollower "D5" "0.64"
Set-TextCell $wsFollower "E5" "1.9"
Set-TextCell $wsFollower "F5" "6.6000000000000005"

# ---------------------------------------------------------------
# Sheet "Punto_modificado" (4th sheet) - x_1, x_2, y_1, y_2 values
# ---------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)

Set-TextCell $wsPunto "A2" "52.150000000000006"
Set-TextCell $wsPunto "B2" "20.25"
Set-TextCell $wsPunto "C2" "104.6"
Set-TextCell $wsPunto "D2" "2.05"

# ---------------------------------------------------------------
# Sheet "Vector_bf" (5th sheet)
# ---------------------------------------------------------------
$wsBf = $wb.Worksheets.Item(5)

Set-TextCell $wsBf "A2" "1.8666746608367828"
Set-TextCell $wsBf "A3" "-2.5326554748285037"

# ---------------------------------------------------------------
# Sheet "Vector_BF" (6th sheet)
# ---------------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)

Set-TextCell $wsBF "A2" "10.200000000000001"
Set-TextCell $wsBF "A3" "-5.800000000000001"
Set-TextCell $wsBF "A4" "-10.79781576900275"
Set-TextCell $wsBF "A5" "-9.156942076603846"

# ---------------------------------------------------------------
# Sheet "Vector_Alpha" (7th sheet) - plain numeric cells (not text)
# ---------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item(7)

$wsAlpha.Range("A2").Value = 0.75
$wsAlpha.Range("A3").Value = 1.0499999999999998
